# Fixed some empty cells in WC file header listing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing D/E/F values for rows 2,3,4,5,9,13 ---
# (these rows were missing no_longer_used / dap_version / userguide_version
# values that the rest of the sheet already has)
$rows = @(2,3,4,5,9,13)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = 0
    $ws.Range("E$r").Value = "3.0.501.0"
    $ws.Range("F$r").Value = 201903
}

# --- Adjust column widths so the newly-visible data is readable ---
# (target stored widths, in "characters": 11.7109375, 15.5703125, 17.28515625,
#  13.85546875, 19.42578125, 22.28515625, 24.42578125 - the ColumnWidth values
#  below are chosen so the engine's internal rounding lands on the closest
#  achievable stored width)
$ws.Columns.Item(1).ColumnWidth = 10.833333333333334
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666
$ws.Columns.Item(4).ColumnWidth = 16.5
$ws.Columns.Item(5).ColumnWidth = 13.0
$ws.Columns.Item(6).ColumnWidth = 18.666666666666668
$ws.Columns.Item(7).ColumnWidth = 21.5
$ws.Columns.Item(8).ColumnWidth = 23.666666666666668

# --- Turn the data range into a formatted Excel table ---
$range = $ws.Range("A1:H95")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight8"

# --- Select A2:H2 to match the saved selection state ---
$ws.Range("A2:H2").Select() | Out-Null

$wb.Save()
